# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'27.920.13"
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').Value = "'1.811.32"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.63%  '
$ws.Range('D4').Value = "'1.002"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.28%  '
$ws.Range('D5').Value = "'309.81"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.69%  '
$ws.Range('E6').Value = '  -0.28%  '
$ws.Range('D7').Value = "'0.4946"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.26%  '
$ws.Range('D8').Value = "'0.3870"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.94%  '
$ws.Range('D9').Value = "'0.09799"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +25.98%  '
$ws.Range('D10').Value = "'1.101"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.55%  '
$ws.Range('D11').Value = "'40.94"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.67%  '
$ws.Range('D12').Value = "'6.450"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.33%  '
$ws.Range('D13').Value = "'20.54"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.94%  '
$ws.Range('D14').Value = "'1.001"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.33%  '
$ws.Range('D15').Value = "'1.811.73"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.61%  '
$ws.Range('D16').Value = "'7.296"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.80%  '
$ws.Range('D17').Value = "'0.00001138"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +6.40%  '
$ws.Range('D18').Value = "'92.67"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.02%  '
$ws.Range('D19').Value = "'0.06612"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.30%  '
$ws.Range('E20').Value = '  -0.29%  '
$ws.Range('D21').Value = "'17.09"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.58%  '
$ws.Range('D22').Value = "'5.931"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.39%  '
$ws.Range('D23').Value = "'27.972.23"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.20%  '
$ws.Range('E24').Value = '  +2.17%  '
$ws.Range('D25').Value = "'2.247"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('D26').Value = "'158.74"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.30%  '
$ws.Range('D27').Value = "'2.021.99"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.78%  '
$ws.Range('D28').Value = "'20.62"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.11%  '
$ws.Range('D29').Value = "'2.396"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.05%  '
$ws.Range('D30').Value = "'127.24"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.09%  '
$ws.Range('D31').Value = "'0.1058"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.86%  '
$ws.Range('E32').Value = '  +0.49%  '
$ws.Range('D33').Value = "'5.574"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.87%  '
$ws.Range('D34').Value = "'3.619"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.24%  '
$ws.Range('D35').Value = "'0.06749"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.73%  '
$ws.Range('D36').Value = "'8.982"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.56%  '
$ws.Range('D37').Value = "'0.02327"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.06%  '
$ws.Range('D38').Value = "'0.2138"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.75%  '
$ws.Range('D39').Value = "'4.938"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.23%  '
$ws.Range('D40').Value = "'11.24"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.88%  '
$ws.Range('D41').Value = "'0.6194"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.81%  '
$ws.Range('E42').Value = '  -0.26%  '
$ws.Range('D43').Value = "'1.144"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.42%  '
$ws.Range('E44').Value = '  +0.88%  '
$ws.Range('D45').Value = "'0.5878"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.75%  '
$ws.Range('D46').Value = "'3.689"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.81%  '
$ws.Range('D47').Value = "'1.277"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.75%  '
$ws.Range('D48').Value = "'122.36"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.06%  '
$ws.Range('D49').Value = "'1.932"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.21%  '
$ws.Range('D50').Value = "'1.176"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.89%  '
$ws.Range('D51').Value = "'0.06788"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.06%  '
